$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------------

# Turn a paragraph's whole text into "<boldText><suffixText>" where boldText
# gets the built-in "Strong" character style (matches <w:rStyle w:val="8"/>)
# and the suffix keeps plain paragraph-level formatting.
function Set-BoldPrefixParagraph($paraIndex, $boldText, $suffixText) {
    $p = $d.Paragraphs.Item($paraIndex)
    $rng = $p.Range
    $rng.MoveEnd(1, -1) | Out-Null
    $startPos = $rng.Start
    $rng.Text = $boldText + $suffixText
    if ($boldText.Length -gt 0) {
        $boldRng = $d.Range($startPos, $startPos + $boldText.Length)
        $boldRng.Style = "Strong"
    }
}

# Build a paragraph made of alternating plain/bold/plain/... runs.
# $parts is an array of @{Text="..."; Bold=$true/$false}
function Set-MixedParagraph($paraIndex, $parts) {
    $p = $d.Paragraphs.Item($paraIndex)
    $rng = $p.Range
    $rng.MoveEnd(1, -1) | Out-Null
    $startPos = $rng.Start
    $full = ""
    foreach ($part in $parts) { $full = $full + $part.Text }
    $rng.Text = $full
    $pos = $startPos
    foreach ($part in $parts) {
        $len = $part.Text.Length
        if ($part.Bold -and $len -gt 0) {
            $sub = $d.Range($pos, $pos + $len)
            $sub.Style = "Strong"
        }
        $pos = $pos + $len
    }
}

# Insert a brand-new bullet paragraph right after paragraph #paraIndex,
# cloning its paragraph style but with a deeper (1440 dxa / 72 pt) left
# indent, and filled with $text.
function Add-Bullet($paraIndex, $text) {
    $p = $d.Paragraphs.Item($paraIndex)
    $rng = $p.Range
    $rng.MoveEnd(1, -1) | Out-Null
    $rng.InsertParagraphAfter() | Out-Null
    $newPara = $d.Paragraphs.Item($paraIndex + 1)
    $newPara.Range.ParagraphFormat.LeftIndent = 72
    $newRng = $newPara.Range
    $newRng.MoveEnd(1, -1) | Out-Null
    $newRng.Text = $text
    return $newPara
}

# ---------------------------------------------------------------------------
# "Integration with ML Pipelines:" bullet list (paragraphs 25-27 originally)
# ---------------------------------------------------------------------------

# Paragraph 25: "Feature generation in Snowflake -> stored in tables"
Set-BoldPrefixParagraph 25 "External Model Training" ":"
Add-Bullet 25 "Extract features using SQL, export to tools like Python or R." | Out-Null
Add-Bullet 26 "Train ML models using scikit-learn, XGBoost, TensorFlow, etc." | Out-Null

# Paragraph is now at index 28 (25 + 3 new total lines incl. itself unchanged count... recompute)
# Original para 26 "Python SDK (or Snowpark) ..." shifted by +2 (two bullets inserted above)
Set-BoldPrefixParagraph 28 "Python + Snowpark" ":"
Add-Bullet 28 "placeholder" | Out-Null
Set-MixedParagraph 29 @(
    @{Text="Use "; Bold=$false},
    @{Text="Snowpark for Python"; Bold=$true},
    @{Text=" to write and execute Python ML code inside Snowflake."; Bold=$false}
)
Add-Bullet 29 "Allows preprocessing, model inference, and even light training directly in Snowflake." | Out-Null

# Original para 27 "Features fetched into memory for training" shifted by +4 total new paragraphs so far
Set-BoldPrefixParagraph 32 "Integration with ML Platforms" ":"
Add-Bullet 32 "Seamless connections with AWS SageMaker, Databricks, and others via connectors." | Out-Null
Add-Bullet 33 "Can schedule ETL jobs and serve features to online models." | Out-Null

Write-Output "Bullets done. p35=[$($d.Paragraphs.Item(35).Range.Text)]"
